$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 3998.75
$ws.Range("J10").Value = 3998.75
$ws.Range("L10").Value = 3998.75
$ws.Range("N10").Value = -4584.75

$ws.Range("H112").Value = 921.6
$ws.Range("J112").Value = 901.9048
$ws.Range("L112").Value = 2705.7144
$ws.Range("N112").Value = -4921.7144

$ws.Range("H116").Value = 216666.67
$ws.Range("I116").Value = 550000
$ws.Range("J116").Value = 50000
$ws.Range("K116").Value = 550000
$ws.Range("L116").Value = 50000
$ws.Range("M116").Value = -546558
$ws.Range("N116").Value = -56884

$ws.Range("H132").Value = 2752.383
$ws.Range("I132").Value = 1150.8286
$ws.Range("K132").Value = 3452.4858
$ws.Range("M132").Value = -922.4858000000004

$ws.Range("H137").Value = 55559080
$ws.Range("I137").Value = 76926104
$ws.Range("J137").Value = 4810
$ws.Range("K137").Value = 230778312
$ws.Range("L137").Value = 14430
$ws.Range("M137").Value = -230775762
$ws.Range("N137").Value = -19530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1517
$ws.Range("I61").Value = 1360.841
$ws.Range("K61").Value = 1360.841
$ws.Range("M61").Value = -1148.841

$ws.Range("H74").Value = 2094.3809
$ws.Range("I74").Value = 971.75
$ws.Range("J74").Value = 5686.8
$ws.Range("K74").Value = 971.75
$ws.Range("L74").Value = 5686.8
$ws.Range("M74").Value = -97.75
$ws.Range("N74").Value = -7434.8

$ws.Range("H77").Value = 2094.3809
$ws.Range("I77").Value = 971.75
$ws.Range("J77").Value = 5686.8
$ws.Range("K77").Value = 4858.75
$ws.Range("L77").Value = 28434
$ws.Range("M77").Value = -490.75
$ws.Range("N77").Value = -37170

$ws.Range("H81").Value = 74287.164
$ws.Range("J81").Value = 74287.164
$ws.Range("L81").Value = 74287.164
$ws.Range("N81").Value = -76283.164

$ws.Range("H84").Value = 74287.164
$ws.Range("J84").Value = 74287.164
$ws.Range("L84").Value = 222861.492
$ws.Range("N84").Value = -232845.492

$ws.Range("H118").Value = 44000
$ws.Range("J118").Value = 44000
$ws.Range("L118").Value = 44000
$ws.Range("N118").Value = -47314

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1517
$ws.Range("I136").Value = 1360.841
$ws.Range("K136").Value = 4082.523
$ws.Range("M136").Value = -1532.523

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2385.8
$ws.Range("I134").Value = 1433.5319
$ws.Range("K134").Value = 4300.5957
$ws.Range("M134").Value = -1765.5957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6456.3945
$ws.Range("I31").Value = 6347.3687
$ws.Range("K31").Value = 6347.3687
$ws.Range("M31").Value = -6052.3687

$ws.Range("H34").Value = 6456.3945
$ws.Range("I34").Value = 6347.3687
$ws.Range("K34").Value = 6347.3687
$ws.Range("M34").Value = -6145.3687

$ws.Range("H58").Value = 1578.1842
$ws.Range("I58").Value = 1165.4722
$ws.Range("J58").Value = 9007
$ws.Range("K58").Value = 1165.4722
$ws.Range("L58").Value = 9007
$ws.Range("M58").Value = -962.4721999999999
$ws.Range("N58").Value = -9413

$ws.Range("H105").Value = 1302.28
$ws.Range("I105").Value = 993.93335
$ws.Range("J105").Value = 1764.8
$ws.Range("K105").Value = 993.93335
$ws.Range("L105").Value = 1764.8
$ws.Range("M105").Value = 753.06665
$ws.Range("N105").Value = -5258.8

$ws.Range("H132").Value = 142858540
$ws.Range("I132").Value = 142858540
$ws.Range("K132").Value = 428575620
$ws.Range("M132").Value = -428573090

$ws.Range("H136").Value = 1578.1842
$ws.Range("I136").Value = 1165.4722
$ws.Range("J136").Value = 9007
$ws.Range("K136").Value = 3496.4166
$ws.Range("L136").Value = 27021
$ws.Range("M136").Value = -946.4165999999996
$ws.Range("N136").Value = -32121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2500
$ws.Range("J62").Value = 2500
$ws.Range("L62").Value = 7500
$ws.Range("N62").Value = -8872

$ws.Range("H65").Value = 2500
$ws.Range("J65").Value = 2500
$ws.Range("L65").Value = 22500
$ws.Range("N65").Value = -29364

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1087.8572
$ws.Range("I113").Value = 903
$ws.Range("K113").Value = 903
$ws.Range("M113").Value = 1267

$ws.Range("H126").Value = 2625.0625
$ws.Range("I126").Value = 2503.3635
$ws.Range("K126").Value = 7510.0905
$ws.Range("M126").Value = -5040.0905

$ws.Range("H132").Value = 200052000
$ws.Range("I132").Value = 250027500
$ws.Range("K132").Value = 750082500
$ws.Range("M132").Value = -750079970

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1070
$ws.Range("I46").Value = 914.2857
$ws.Range("J46").Value = 1433.3334
$ws.Range("K46").Value = 914.2857
$ws.Range("L46").Value = 1433.3334
$ws.Range("M46").Value = -726.2857
$ws.Range("N46").Value = -1809.3334

$ws.Range("H132").Value = 3124.625
$ws.Range("I132").Value = 3142.4285
$ws.Range("K132").Value = 9427.2855
$ws.Range("M132").Value = -6897.2855

$ws.Range("H136").Value = 3545.2
$ws.Range("I136").Value = 2244.4614
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 6733.3842
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -4183.3842
$ws.Range("N136").Value = -41100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1480.1052
$ws.Range("I122").Value = 1294.8235
$ws.Range("J122").Value = 3055
$ws.Range("K122").Value = 3884.4705
$ws.Range("L122").Value = 9165
$ws.Range("M122").Value = -1434.4705
$ws.Range("N122").Value = -14065

$ws.Range("H126").Value = 2282
$ws.Range("I126").Value = 1709.0588
$ws.Range("J126").Value = 3093.6667
$ws.Range("K126").Value = 5127.1764
$ws.Range("L126").Value = 9281.000100000001
$ws.Range("M126").Value = -2657.1764
$ws.Range("N126").Value = -14221.0001

$ws.Range("H132").Value = 7411403.5
$ws.Range("I132").Value = 11113928
$ws.Range("J132").Value = 6355.3335
$ws.Range("K132").Value = 33341784
$ws.Range("L132").Value = 19066.0005
$ws.Range("M132").Value = -33339254
